$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 64) with the next date in the series and the
# carried-forward price values (same as the prior row, 63).
# Force the date cell to be stored as literal text (matching the existing
# "Date" column cells, which are plain strings rather than Excel date
# serials), then restore the cell to the default "Normal" style so no
# lingering number-format override is left on the cell.
$ws.Range("A64").NumberFormat = "@"
$ws.Range("A64").Value = "2025-10-18"
$ws.Range("A64").Style = "Normal"

$ws.Range("B64").Value = 52.91999816894531
$ws.Range("C64").Value = 396.6000061035156
$ws.Range("D64").Value = 342.6499938964844
